$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to Text format first so Excel does not auto-convert
# numeric-looking price strings (e.g. "519.83") into floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "58.234.46"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").Value = "2.610.05"
$ws.Range("E3").Value = "  -4.22%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "519.83"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "142.90"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").Value = "6.67"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "3.067.53"
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("D14").Value = "58.213.78"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "20.96"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "2.620.79"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").Value = "4.40"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").Value = "335.50"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "10.37"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "64.57"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "0.0₃0791"
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").Value = "6.62"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D31").Value = "1.59"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "18.76"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").Value = "150.10"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("D36").Value = "0.888"
$ws.Range("E36").Value = "  -5.59%  "
$ws.Range("D37").Value = "0.853"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").Value = "36.27"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "0.599"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("D44").Value = "269.57"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D46").Value = "19.17"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "2.035.98"
$ws.Range("E48").Value = "  -4.90%  "
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "4.69"
$ws.Range("E50").Value = "  -4.85%  "
$ws.Range("E51").Value = "  -4.33%  "
